$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for existing rows 2-6 (B, C, E, F columns; A and D unchanged)
$rows = @(
    @{ r=2;  A=0; B=1;  C=458.8499701346049; D=1500; E=457.524951905785;  F=809.2470893547023 },
    @{ r=3;  A=1; B=1;  C=566.9923819505786; D=1500; E=544.4008421964322; F=935.9457710565886 },
    @{ r=4;  A=2; B=1;  C=3897.301207536528; D=1500; E=421.3954489148648; F=822.6894271705743 },
    @{ r=5;  A=3; B=1;  C=583.1179766099127; D=1500; E=541.0500768672667; F=1010.930492410078 },
    @{ r=6;  A=4; B=1;  C=507.2009852429924; D=1500; E=466.9654057349896; F=886.1562661123647 },
    @{ r=7;  A=5; B=64; C=2135.527146272234; D=1500; E=527.5487730867658; F=950.8942286852177 },
    @{ r=8;  A=6; B=1;  C=547.8085094871341; D=1500; E=500.9139139457083; F=935.4559132762812 },
    @{ r=9;  A=7; B=1;  C=557.6323942244364; D=1500; E=511.2995171506442; F=897.3191306581153 },
    @{ r=10; A=8; B=1;  C=491.9717572709463; D=1500; E=469.7719128676806; F=867.534308455666  },
    @{ r=11; A=9; B=1;  C=498.5914671003533; D=1500; E=429.929044104556;  F=847.4310521269317 }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
}

# Copy the style from A6 (existing styled cell) to the new A7:A11 cells
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A7:A11").PasteSpecial(-4122) | Out-Null
